$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.181.38"
$ws.Range("E2").Value = "  +0.49%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.835.41"
$ws.Range("E3").Value = "  +0.19%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9987"
$ws.Range("E4").Value = "  -0.03%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "241.97"
$ws.Range("E5").Value = "  +1.17%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6572"
$ws.Range("E6").Value = "  -1.81%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.9994"
$ws.Range("E7").Value = "  -0.06%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07438"
$ws.Range("E8").Value = "  +0.35%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2925"
$ws.Range("E9").Value = "  -0.90%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "22.97"
$ws.Range("E10").Value = "  +1.02%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07737"
$ws.Range("E11").Value = "  +1.26%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.823.71"
$ws.Range("E12").Value = "  -0.60%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.979"
$ws.Range("E13").Value = "  -0.52%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6657"
$ws.Range("E14").Value = "  -1.03%  "

$ws.Range("E15").Value = "  -3.99%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.121"
$ws.Range("E16").Value = "  -0.25%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008497"
$ws.Range("E17").Value = "  +3.24%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "29.143.51"
$ws.Range("E18").Value = "  +0.39%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.984.48"
$ws.Range("E19").Value = "  -4.41%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "227.34"
$ws.Range("E20").Value = "  -0.23%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.46"
$ws.Range("E21").Value = "  +0.20%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.9992"
$ws.Range("E22").Value = "  -0.08%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.138"
$ws.Range("E23").Value = "  -2.02%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.9994"
$ws.Range("E24").Value = "  -0.01%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "158.56"

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.594"
$ws.Range("E26").Value = "  -0.75%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.1394"
$ws.Range("E27").Value = "  -2.56%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "17.93"
$ws.Range("E28").Value = "  -0.21%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.518"

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.115"
$ws.Range("E30").Value = "  -2.76%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.046"

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.193"
$ws.Range("E32").Value = "  -0.32%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05257"
$ws.Range("E33").Value = "  -2.19%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.866"
$ws.Range("E34").Value = "  +0.81%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7336"
$ws.Range("E35").Value = "  -2.16%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.144"
$ws.Range("E36").Value = "  +1.53%  "

$ws.Range("E37").Value = "  -1.15%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.304.67"
$ws.Range("E38").Value = "  +0.68%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01797"
$ws.Range("E39").Value = "  -0.61%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.732"
$ws.Range("E40").Value = "  +1.14%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9210"
$ws.Range("E41").Value = "  +0.01%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.08778"
$ws.Range("E42").Value = "  +10.08%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.966"
$ws.Range("E43").Value = "  -0.84%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.9988"
$ws.Range("E44").Value = "  +0.05%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "102.23"
$ws.Range("E45").Value = "  -1.77%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.964.94"
$ws.Range("E46").Value = "  -0.57%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5140"
$ws.Range("E47").Value = "  -0.68%  "

$ws.Range("B48").Value = "RenderToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.758"
$ws.Range("E48").Value = "  +0.37%  "

$ws.Range("B49").Value = "Aave"
$ws.Range("C49").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "63.40"
$ws.Range("E49").Value = "  -0.31%  "

$ws.Range("B50").Value = "BabyDogeCoin"
$ws.Range("C50").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00000000120"
$ws.Range("E50").Value = "  -2.58%  "

$ws.Range("E51").Value = "  -1.12%  "
